$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DecemberRaw")

$ws.Cells.Item(1, 1).Value = "Library"
$ws.Cells.Item(1, 2).Value = "Items owned by this library checked out at this library this month"
$ws.Cells.Item(1, 3).Value = "Items owned by other libraries checked out at this library this month"
$ws.Cells.Item(1, 4).Value = "Total circulation this month"

$ws.Cells.Item(2, 1).Value = "Atchison Public Library"
$ws.Cells.Item(2, 2).Value = 3335
$ws.Cells.Item(2, 3).Value = 1324
$ws.Cells.Item(2, 4).Value = 4659
$ws.Cells.Item(3, 1).Value = "Baldwin City Public Library"
$ws.Cells.Item(3, 2).Value = 2200
$ws.Cells.Item(3, 3).Value = 510
$ws.Cells.Item(3, 4).Value = 2710
$ws.Cells.Item(4, 1).Value = "Basehor Community Library"
$ws.Cells.Item(4, 2).Value = 6997
$ws.Cells.Item(4, 3).Value = 975
$ws.Cells.Item(4, 4).Value = 7972
$ws.Cells.Item(5, 1).Value = "Bern Community Library"
$ws.Cells.Item(5, 2).Value = 93
$ws.Cells.Item(5, 3).Value = 29
$ws.Cells.Item(5, 4).Value = 122
$ws.Cells.Item(6, 1).Value = "Bonner Springs City Library"
$ws.Cells.Item(6, 2).Value = 4206
$ws.Cells.Item(6, 3).Value = 836
$ws.Cells.Item(6, 4).Value = 5042
$ws.Cells.Item(7, 1).Value = "Burlingame Community Library"
$ws.Cells.Item(7, 2).Value = 353
$ws.Cells.Item(7, 3).Value = 154
$ws.Cells.Item(7, 4).Value = 507
$ws.Cells.Item(8, 1).Value = "Carbondale City Library"
$ws.Cells.Item(8, 2).Value = 521
$ws.Cells.Item(8, 3).Value = 140
$ws.Cells.Item(8, 4).Value = 661
$ws.Cells.Item(9, 1).Value = "Centralia Community Library"
$ws.Cells.Item(9, 2).Value = 202
$ws.Cells.Item(9, 3).Value = 32
$ws.Cells.Item(9, 4).Value = 234
$ws.Cells.Item(10, 1).Value = "Corning City Library"
$ws.Cells.Item(10, 2).Value = 43
$ws.Cells.Item(10, 3).Value = 3
$ws.Cells.Item(10, 4).Value = 46
$ws.Cells.Item(11, 1).Value = "Digital Content"
$ws.Cells.Item(12, 1).Value = "Doniphan County Library - Elwood"
$ws.Cells.Item(12, 2).Value = 77
$ws.Cells.Item(12, 3).Value = 19
$ws.Cells.Item(12, 4).Value = 96
$ws.Cells.Item(13, 1).Value = "Doniphan County Library - Highland"
$ws.Cells.Item(13, 2).Value = 156
$ws.Cells.Item(13, 3).Value = 91
$ws.Cells.Item(13, 4).Value = 247
$ws.Cells.Item(14, 1).Value = "Doniphan County Library - Troy"
$ws.Cells.Item(14, 2).Value = 396
$ws.Cells.Item(14, 3).Value = 144
$ws.Cells.Item(14, 4).Value = 540
$ws.Cells.Item(15, 1).Value = "Doniphan County Library - Wathena"
$ws.Cells.Item(15, 2).Value = 199
$ws.Cells.Item(15, 3).Value = 43
$ws.Cells.Item(15, 4).Value = 242
$ws.Cells.Item(16, 1).Value = "Effingham Community Library"
$ws.Cells.Item(16, 2).Value = 187
$ws.Cells.Item(16, 3).Value = 64
$ws.Cells.Item(16, 4).Value = 251
$ws.Cells.Item(17, 1).Value = "Eudora Community Library"
$ws.Cells.Item(17, 2).Value = 1131
$ws.Cells.Item(17, 3).Value = 499
$ws.Cells.Item(17, 4).Value = 1630
$ws.Cells.Item(18, 1).Value = "Everest, Barnes Reading Room"
$ws.Cells.Item(18, 2).Value = 47
$ws.Cells.Item(18, 3).Value = 64
$ws.Cells.Item(18, 4).Value = 111
$ws.Cells.Item(19, 1).Value = "Hiawatha, Morrill Public Library"
$ws.Cells.Item(19, 2).Value = 1356
$ws.Cells.Item(19, 3).Value = 427
$ws.Cells.Item(19, 4).Value = 1783
$ws.Cells.Item(20, 1).Value = "Highland Community College"
$ws.Cells.Item(20, 2).Value = 10
$ws.Cells.Item(20, 3).Value = 20
$ws.Cells.Item(20, 4).Value = 30
$ws.Cells.Item(21, 1).Value = "Holton, Beck-Bookman Library"
$ws.Cells.Item(21, 2).Value = 1350
$ws.Cells.Item(21, 3).Value = 477
$ws.Cells.Item(21, 4).Value = 1827
$ws.Cells.Item(22, 1).Value = "Horton Public Library"
$ws.Cells.Item(22, 2).Value = 169
$ws.Cells.Item(22, 3).Value = 45
$ws.Cells.Item(22, 4).Value = 214
$ws.Cells.Item(23, 1).Value = "Lansing Community Library"
$ws.Cells.Item(23, 2).Value = 1584
$ws.Cells.Item(23, 3).Value = 467
$ws.Cells.Item(23, 4).Value = 2051
$ws.Cells.Item(24, 1).Value = "Leavenworth Public Library"
$ws.Cells.Item(24, 2).Value = 7047
$ws.Cells.Item(24, 3).Value = 1420
$ws.Cells.Item(24, 4).Value = 8467
$ws.Cells.Item(25, 1).Value = "Linwood Community Library"
$ws.Cells.Item(25, 2).Value = 516
$ws.Cells.Item(25, 3).Value = 177
$ws.Cells.Item(25, 4).Value = 693
$ws.Cells.Item(26, 1).Value = "Louisburg Library"
$ws.Cells.Item(27, 1).Value = "Lyndon Carnegie Library"
$ws.Cells.Item(27, 2).Value = 351
$ws.Cells.Item(27, 3).Value = 227
$ws.Cells.Item(27, 4).Value = 578
$ws.Cells.Item(28, 1).Value = "McLouth Public Library"
$ws.Cells.Item(28, 2).Value = 216
$ws.Cells.Item(28, 3).Value = 73
$ws.Cells.Item(28, 4).Value = 289
$ws.Cells.Item(29, 1).Value = "Meriden-Ozawkie Public Library"
$ws.Cells.Item(29, 2).Value = 1080
$ws.Cells.Item(29, 3).Value = 539
$ws.Cells.Item(29, 4).Value = 1619
$ws.Cells.Item(30, 1).Value = "Northeast Kansas Library System"
$ws.Cells.Item(30, 2).Value = 11
$ws.Cells.Item(30, 3).Value = 55
$ws.Cells.Item(30, 4).Value = 66
$ws.Cells.Item(31, 1).Value = "Nortonville Public Library"
$ws.Cells.Item(31, 2).Value = 234
$ws.Cells.Item(31, 3).Value = 47
$ws.Cells.Item(31, 4).Value = 281
$ws.Cells.Item(32, 1).Value = "Osage City Library"
$ws.Cells.Item(32, 2).Value = 1147
$ws.Cells.Item(32, 3).Value = 354
$ws.Cells.Item(32, 4).Value = 1501
$ws.Cells.Item(33, 1).Value = "Osawatomie Public Library"
$ws.Cells.Item(33, 2).Value = 942
$ws.Cells.Item(33, 3).Value = 322
$ws.Cells.Item(33, 4).Value = 1264
$ws.Cells.Item(34, 1).Value = "Oskaloosa Public Library"
$ws.Cells.Item(34, 2).Value = 446
$ws.Cells.Item(34, 3).Value = 178
$ws.Cells.Item(34, 4).Value = 624
$ws.Cells.Item(35, 1).Value = "Ottawa Library"
$ws.Cells.Item(35, 2).Value = 4800
$ws.Cells.Item(35, 3).Value = 732
$ws.Cells.Item(35, 4).Value = 5532
$ws.Cells.Item(36, 1).Value = "Overbrook Public Library"
$ws.Cells.Item(36, 2).Value = 664
$ws.Cells.Item(36, 3).Value = 123
$ws.Cells.Item(36, 4).Value = 787
$ws.Cells.Item(37, 1).Value = "Paola Free Library"
$ws.Cells.Item(37, 2).Value = 2777
$ws.Cells.Item(37, 3).Value = 484
$ws.Cells.Item(37, 4).Value = 3261
$ws.Cells.Item(38, 1).Value = "Perry-Lecompton Community Library"
$ws.Cells.Item(38, 2).Value = 97
$ws.Cells.Item(38, 3).Value = 31
$ws.Cells.Item(38, 4).Value = 128
$ws.Cells.Item(39, 1).Value = "Pomona Community Library"
$ws.Cells.Item(39, 2).Value = 61
$ws.Cells.Item(39, 3).Value = 74
$ws.Cells.Item(39, 4).Value = 135
$ws.Cells.Item(40, 1).Value = "Prairie Hills Schools - Axtell Public School"
$ws.Cells.Item(40, 2).Value = 267
$ws.Cells.Item(40, 3).Value = 21
$ws.Cells.Item(40, 4).Value = 288
$ws.Cells.Item(41, 1).Value = "Prairie Hills Schools - Sabetha Elementary School"
$ws.Cells.Item(41, 2).Value = 1705
$ws.Cells.Item(41, 3).Value = 88
$ws.Cells.Item(41, 4).Value = 1793
$ws.Cells.Item(42, 1).Value = "Prairie Hills Schools - Sabetha High School"
$ws.Cells.Item(42, 2).Value = 23
$ws.Cells.Item(42, 3).Value = 3
$ws.Cells.Item(42, 4).Value = 26
$ws.Cells.Item(43, 1).Value = "Prairie Hills Schools - Sabetha Middle School"
$ws.Cells.Item(43, 2).Value = 105
$ws.Cells.Item(43, 3).Value = 6
$ws.Cells.Item(43, 4).Value = 111
$ws.Cells.Item(44, 1).Value = "Prairie Hills Schools - Wetmore Academic Center (Permanently closed)"
$ws.Cells.Item(45, 1).Value = "Richmond Public Library"
$ws.Cells.Item(45, 2).Value = 281
$ws.Cells.Item(45, 3).Value = 63
$ws.Cells.Item(45, 4).Value = 344
$ws.Cells.Item(46, 1).Value = "Rossville Community Library"
$ws.Cells.Item(46, 2).Value = 1018
$ws.Cells.Item(46, 3).Value = 399
$ws.Cells.Item(46, 4).Value = 1417
$ws.Cells.Item(47, 1).Value = "Sabetha, Mary Cotton Library"
$ws.Cells.Item(47, 2).Value = 2395
$ws.Cells.Item(47, 3).Value = 792
$ws.Cells.Item(47, 4).Value = 3187
$ws.Cells.Item(48, 1).Value = "Seneca Free Library"
$ws.Cells.Item(48, 2).Value = 1360
$ws.Cells.Item(48, 3).Value = 155
$ws.Cells.Item(48, 4).Value = 1515
$ws.Cells.Item(49, 1).Value = "Silver Lake Library"
$ws.Cells.Item(49, 2).Value = 780
$ws.Cells.Item(49, 3).Value = 650
$ws.Cells.Item(49, 4).Value = 1430
$ws.Cells.Item(50, 1).Value = "Tonganoxie Public Library"
$ws.Cells.Item(50, 2).Value = 2671
$ws.Cells.Item(50, 3).Value = 777
$ws.Cells.Item(50, 4).Value = 3448
$ws.Cells.Item(51, 1).Value = "Valley Falls, Delaware Township Library"
$ws.Cells.Item(51, 2).Value = 371
$ws.Cells.Item(51, 3).Value = 161
$ws.Cells.Item(51, 4).Value = 532
$ws.Cells.Item(52, 1).Value = "Wellsville City Library"
$ws.Cells.Item(52, 2).Value = 599
$ws.Cells.Item(52, 3).Value = 209
$ws.Cells.Item(52, 4).Value = 808
$ws.Cells.Item(53, 1).Value = "Wetmore Public Library"
$ws.Cells.Item(53, 2).Value = 133
$ws.Cells.Item(53, 3).Value = 132
$ws.Cells.Item(53, 4).Value = 265
$ws.Cells.Item(54, 1).Value = "Williamsburg Community Library"
$ws.Cells.Item(54, 2).Value = 184
$ws.Cells.Item(54, 3).Value = 30
$ws.Cells.Item(54, 4).Value = 214
$ws.Cells.Item(55, 1).Value = "Winchester Public Library"
$ws.Cells.Item(55, 2).Value = 287
$ws.Cells.Item(55, 3).Value = 519
$ws.Cells.Item(55, 4).Value = 806

$excel.CalculateFull()
